$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store its value as TEXT (matches the source data,
# which keeps numeric-looking strings like ids/prices as plain text), while
# leaving the cell's style untouched (reset to "Normal" after the write so no
# stray number-format style sticks around).
function Set-TextCell($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

# --- Row 21 becomes what used to be row 22 ("Tela Viva ...") ---
Set-TextCell $ws.Cells.Item(21, 1) "6996030"
Set-TextCell $ws.Cells.Item(21, 2) "Tela Viva Haushaltspapier 3-lagig 4 Rollen"
Set-TextCell $ws.Cells.Item(21, 3) "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/tela-viva-haushaltspapier-3-lagig-4-rollen/p/6996030"
Set-TextCell $ws.Cells.Item(21, 4) "200BLT"
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(21, 6).Value2 = 0
Set-TextCell $ws.Cells.Item(21, 7) "Tela"
Set-TextCell $ws.Cells.Item(21, 8) "5.95"
$ws.Cells.Item(21, 9).ClearContents()
$ws.Cells.Item(21, 10).ClearContents()
$ws.Cells.Item(21, 11).ClearContents()
$ws.Cells.Item(21, 12).ClearContents()
Set-TextCell $ws.Cells.Item(21, 13) "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
Set-TextCell $ws.Cells.Item(21, 14) "Tela Viva Haushaltspapier 3-lagig 4 Rollen 5.95 Schweizer Franken"

# --- Row 22 becomes what used to be row 21 ("Tempo Bamboo Eco") ---
Set-TextCell $ws.Cells.Item(22, 1) "6868354"
Set-TextCell $ws.Cells.Item(22, 2) "Tempo Bamboo Eco"
Set-TextCell $ws.Cells.Item(22, 3) "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco/p/6868354"
Set-TextCell $ws.Cells.Item(22, 4) "12ST"
$ws.Cells.Item(22, 5).Value2 = 1
$ws.Cells.Item(22, 6).Value2 = 3
Set-TextCell $ws.Cells.Item(22, 7) "Tempo"
Set-TextCell $ws.Cells.Item(22, 8) "3.95"
Set-TextCell $ws.Cells.Item(22, 9) "0.33/1ST"
Set-TextCell $ws.Cells.Item(22, 10) "Preis pro 1 Stück"
Set-TextCell $ws.Cells.Item(22, 11) "0.33"
Set-TextCell $ws.Cells.Item(22, 12) "1ST"
Set-TextCell $ws.Cells.Item(22, 13) "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
Set-TextCell $ws.Cells.Item(22, 14) "Tempo Bamboo Eco 3.95 Schweizer Franken"

# --- Refresh the scrape timestamp (column O) for every data row (2-26) ---
$newTimestamp = "2022-08-26 20:59:18"
for ($r = 2; $r -le 26; $r++) {
    Set-TextCell $ws.Cells.Item($r, 15) $newTimestamp
}
